## Drybar commit by Varaprasad 1-03-2021
## Applies the DryBarTestData.xlsx edits:
##  - shortens the sample product name in V3
##  - fixes up the comma separated header-name list in W11
##  - adds seven new "product name" validation columns (AV:BB) with header +
##    sample row
##  - stores the gift-card code / pin test values as text ("0123") instead of
##    numbers
##  - updates the active selection

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Shorten the product name used for the HydroFlask/DryBar sample row.
$ws.Range("V3").Value = "Liquid Glass"

# 2. Correct wording in the pipe-style header name list on row 11.
$ws.Range("W11").Value = "Hair Products,Hair Tools,Benefits,Gifts & Kits,New,How To & Inspo"

# 3. Gift card code / pin test values become text "0123" (quote-prefixed)
#    instead of the numeric 123.
$ws.Range("X12").Value = "'0123"
$ws.Range("Y12").Value = "'0123"

# 4. New product-name test columns AV:BB with a header (row 1) and a sample
#    value (row 3). Match the existing header look (yellow fill, left/top
#    aligned) used by the rest of row 1.
$headers = @(
    @{ Col = "AV"; Header = "Threedigitproductname"; Sample = "Hai" },
    @{ Col = "AW"; Header = "Fourdigitproductname";  Sample = "Hair" },
    @{ Col = "AX"; Header = "Nameofproduct";         Sample = "Hair Products" },
    @{ Col = "AY"; Header = "dublicateproductname";  Sample = "color care" },
    @{ Col = "AZ"; Header = "invalidname";           Sample = "ppp" },
    @{ Col = "BA"; Header = "Pname";                 Sample = "Shampoos" },
    @{ Col = "BB"; Header = "vicksproductname";      Sample = "humidifiers" }
)

$headerRange = $ws.Range("AV1:BB1")
$headerRange.Interior.Color = 65535
$headerRange.HorizontalAlignment = -4131
$headerRange.VerticalAlignment = -4160

foreach ($h in $headers) {
    $ws.Range($h.Col + "1").Value = $h.Header
    $ws.Range($h.Col + "3").Value = $h.Sample
}

# 5. Column widths - best effort match of the resized/new columns.
$ws.Columns.Item(18).ColumnWidth = 18.33   # R
$ws.Columns.Item(19).ColumnWidth = 25.67   # S
$ws.Columns.Item(22).ColumnWidth = 14.83   # V
$ws.Columns.Item(48).ColumnWidth = 26.5    # AV
$ws.Columns.Item(49).ColumnWidth = 23.17   # AW
$ws.Columns.Item(50).ColumnWidth = 16.67   # AX
$ws.Columns.Item(51).ColumnWidth = 20.33   # AY
$ws.Columns.Item(52).ColumnWidth = 12.5    # AZ
$ws.Columns.Item(53).ColumnWidth = 9.67    # BA
$ws.Columns.Item(54).ColumnWidth = 16.67   # BB

# 6. Selection moves to D11, with no frozen/scrolled top-left cell.
$ws.Range("D11").Select()
